# Super Serious Space Game - List_of_projects.docx
#
# Commit: "Added ship special cooldown and cost"
#
# Marks the following TODO items as done (bold + strikethrough + red),
# matching the look already used elsewhere in the doc for completed items,
# and moves the "_GoBack" last-edit bookmark from the "Splash doesn't
# play after some battles" paragraph to the very end of the document
# (right after "etc" in the last paragraph), since that's where the
# author's last edit for this commit landed.

$d = $word.ActiveDocument

function Mark-Done($paraIndex, $bold) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Font.StrikeThrough = 1
    $r.Font.Color = 255        # wdColorRed == RGB(255,0,0) -> w:color val="FF0000"
    if ($bold) {
        $r.Font.Bold = 1
    }
}

# --- Battle: ---------------------------------------------------------
# "Bring power supply to battle to spawn ship special"
Mark-Done 11 $true
# "Add starting cool" / " " / "down for ship special" (3 runs, 1 paragraph)
Mark-Done 12 $true

# "Find and fix all of the " / "units" / " " / "retrieveUpgrades" / " references"
Mark-Done 20 $true

# --- Player: -----------------------------------------------------------
# "Create player stats resource caps/miner caps/unit caps etc."
Mark-Done 72 $true
# "Add spawn pod type upgrades cap/cost/types etc."
Mark-Done 73 $true

# --- Move the _GoBack bookmark -----------------------------------------
# Currently sits (collapsed) right after "...t play after some battles".
# It needs to end up (collapsed) right after "etc" in the very last
# paragraph, before the trailing spell-check proofErr marker.
#
# NOTE: this COM host has a positional quirk where a *collapsed* Range
# sitting exactly at (paragraph.Range.End - 1) cannot be used directly
# as a bookmark anchor. We work around it by temporarily inserting a
# placeholder character after the target point (so the target position
# is no longer right at that boundary), anchoring the bookmark there,
# and then deleting the placeholder -- the bookmark keeps its place.

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$afterEtc = $lastPara.Range.End - 1   # right after "etc", before the para mark

$placeholder = $d.Range($afterEtc, $afterEtc)
$placeholder.InsertAfter([char]1)

$anchor = $d.Range($afterEtc, $afterEtc)
$d.Bookmarks.Add("_GoBack", $anchor)

$bm = $d.Bookmarks.Item("_GoBack")
$placeholderRange = $d.Range($bm.End, $bm.End + 1)
$placeholderRange.Delete()
